$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38 (shifts existing rows 38..158 down to 39..159,
# carrying forward the date-column number format from the row above).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Cells.Item(38, 1).Value = 5
$ws.Cells.Item(38, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(38, 3).Value = "Maule"
$ws.Cells.Item(38, 4).Value = 44453
$ws.Cells.Item(38, 5).Value = 7
$ws.Cells.Item(38, 6).Value = 100112003
$ws.Cells.Item(38, 7).Value = "Ajo"
$ws.Cells.Item(38, 8).Value = "Chino"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 200
$ws.Cells.Item(38, 11).Value = 15000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = 15000
$ws.Cells.Item(38, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(38, 15).Value = "China"
$ws.Cells.Item(38, 16).Value = 1500
$ws.Cells.Item(38, 17).Value = 10
$ws.Cells.Item(38, 18).Value = "Hortaliza"
